# Book1.xlsx - rename the "Completed" column to "Status" and fix up its
# last row value, then leave the sheet in the state it was last saved in
# (column L autosized, cursor on L6, portrait page orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Column L was labelled "Completed" - rename it to "Status".
$ws.Range("L1").Value = "Status"

# Row 6's status value was wrong (2) - correct it to 0.
$ws.Range("L6").Value = 0

# Resize column L to fit the new header/content.
$ws.Columns.Item(12).AutoFit() | Out-Null

# Make sure the sheet prints in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection on L6, where the edit was made.
$ws.Range("L6").Select() | Out-Null

$wb.Save()
